$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.738.85'
$ws.Range('E2').Value = '  +4.74%  '
$ws.Range('D3').Value = '2.281.97'
$ws.Range('E3').Value = '  +2.34%  '
$ws.Range('E4').Value = '  +0.25%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '231.22'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.11%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.627'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.42%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '63.04'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +4.40%  '
$ws.Range('E8').Value = '  +0.10%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.422'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +4.12%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0950'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +5.26%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '57.59'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.17%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '26.06'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +12.89%  '
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').Value = '2.625.28'
$ws.Range('E14').Value = '  +2.55%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '15.82'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.69%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.94'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +5.75%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.812'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '2.288.13'
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('D19').Value = '43.710.38'
$ws.Range('E19').Value = '  +4.76%  '
$ws.Range('D20').Value = '0.0₃0962'
$ws.Range('E20').Value = '  +6.41%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '73.07'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.06%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.17'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.65%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '252.24'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.53%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.63'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +11.20%  '
$ws.Range('E25').Value = '  -0.04%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -2.21%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.84'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +1.14%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '171.11'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.97%  '
$ws.Range('E29').Value = '  -2.00%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '20.49'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +2.80%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.44'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.84%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.75'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +3.82%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.122'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.21%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.0689'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +6.37%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '5.08'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.86%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '4.72'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '6.57'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '3.69'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +1.39%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.36'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.96%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.0249'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +3.57%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '11.01'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +27.96%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '8.67'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '4.65'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +6.19%  '
$ws.Range('B45').Value = 'TerraClassic'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.000221'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -7.39%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.22'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0966'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '98.06'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.63%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '17.08'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +2.80%  '
$ws.Range('D50').Value = '1.483.71'
$ws.Range('E50').Value = '  +1.07%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.31'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.51%  '
